$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (compensate for the runtime's systematic +5/6 char offset) ---
$ws.Columns.Item(2).ColumnWidth = 25 - 0.8333333333333334
$ws.Columns.Item(4).ColumnWidth = 40 - 0.8333333333333334

# --- Header row styling (blue fill, white bold font; bold already set) ---
$headerRange = $ws.Range("A1:E1")
$headerRange.Interior.Color = 12419407
$headerRange.Font.Color = 16777215

# --- Row 2 cell values ---
$ws.Range("A2").Value = "26/06/2025"
$ws.Range("B2").Value = "MEGA FRUVER EL SUPEF"
$ws.Range("C2").Value = "TOTAL NO ENCONTRADO"
$ws.Range("C2").NumberFormat = '"$"#,##0.00'
$ws.Range("D2").Value = "'0.59"
$ws.Range("E2").Value = 'Espacio para
; Logo Corporativo
MEGA FRUVER EL SUPEF
REGALON SAS
Le
Dir.: CALLE 8 26 - 17 |
rE ISA ies
Documento de ingreso
PMP AL Yd
er Cle 26/06/2025, 11 35
CMe Cr AO RE eT BC Reems
NaC M TMC Ne Cat alk 1 eke 1a ee
soporte de uso cantable
ety Consumidor Final
C.c i NIT: 222222222222-7
ital
Vendedor: jhon anderson arango
| or ZT)
0.59
4 yt) RATA) PAR,
PAPA AMARIILA GRANEL / 3001
ce
2 ts Pe oPoUeR
LIMON TAHITI / 1017
E pK y
3 mT ALOK) SSK
es San ad
‘. are
C224.)
Ae)
Rend
ray'

# --- Remove auto row-height stamp introduced by the multi-line E2 text ---
$ws.Rows.Item(2).AutoFit()

